# Auto-generated script to apply cryptos.xlsx diff
# (crypto price/volume refresh + FTXToken/FraxShare row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay stored as text (matches source inlineStr cells)
$textCells = @("D5", "D6", "D7", "D10", "D12", "D14", "D15", "D17", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D34", "D40", "D41", "D43", "D45", "D46", "D48", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.095.43"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.052.24"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "249.94"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "0.668"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "59.72"
$ws.Range("E7").Value = "  +7.96%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "0.0794"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "16.10"
$ws.Range("E12").Value = "  +7.07%  "
$ws.Range("D13").Value = "2.351.66"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "0.835"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").Value = "5.80"
$ws.Range("E15").Value = "  +10.19%  "
$ws.Range("D16").Value = "2.051.31"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "18.36"
$ws.Range("E17").Value = "  +29.28%  "
$ws.Range("D18").Value = "37.077.44"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "76.01"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "5.40"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "238.49"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  +11.31%  "
$ws.Range("D26").Value = "169.66"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").Value = "9.42"
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").Value = "20.20"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +7.98%  "
$ws.Range("D31").Value = "4.81"
$ws.Range("E31").Value = "  +5.12%  "
$ws.Range("D32").Value = "0.0634"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").Value = "0.0892"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  +4.20%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  +14.50%  "
$ws.Range("D41").Value = "5.29"
$ws.Range("E41").Value = "  +20.99%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "17.61"
$ws.Range("E43").Value = "  -2.50%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "97.40"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").Value = "1.295.81"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "2.89"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "6.87"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "3.68"
$ws.Range("E50").Value = "  -19.97%  "
$ws.Range("D51").Value = "2.242.64"
$ws.Range("E51").Value = "  -0.12%  "
